$d = $word.ActiveDocument

$pairs = @(
    @("765÷5=153, 0", "274÷6=45, 4"),
    @("605÷7=86, 3", "187÷6=31, 1"),
    @("893÷8=111, 5", "418÷7=59, 5"),
    @("530÷8=66, 2", "700÷8=87, 4"),
    @("602÷8=75, 2", "780÷7=111, 3"),
    @("110÷9=12, 2", "343÷3=114, 1"),
    @("888÷2=444, 0", "336÷7=48, 0"),
    @("285÷9=31, 6", "418÷7=59, 5"),
    @("871÷2=435, 1", "655÷4=163, 3"),
    @("513÷7=73, 2", "290÷5=58, 0"),
    @("832÷4=208, 0", "682÷6=113, 4"),
    @("793÷3=264, 1", "544÷7=77, 5"),
    @("303÷4=75, 3", "590÷3=196, 2"),
    @("990÷8=123, 6", "679÷4=169, 3"),
    @("117÷2=58, 1", "262÷9=29, 1"),
    @("769÷3=256, 1", "821÷3=273, 2"),
    @("643÷9=71, 4", "119÷5=23, 4"),
    @("991÷6=165, 1", "881÷7=125, 6"),
    @("613÷7=87, 4", "652÷8=81, 4"),
    @("100÷7=14, 2", "957÷3=319, 0"),
    @("191÷5=38, 1", "678÷4=169, 2"),
    @("370÷4=92, 2", "395÷6=65, 5"),
    @("869÷4=217, 1", "515÷8=64, 3"),
    @("995÷8=124, 3", "535÷5=107, 0"),
    @("567÷8=70, 7", "625÷3=208, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
